$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the water_depth and well_depth columns
$ws.Range("H1").Value = "water_depth"
$ws.Range("I1").Value = "well_depth"

# Fill in the water_depth values (column H)
$ws.Range("H2").Value = 132
$ws.Range("H3").Value = 261
$ws.Range("H4").Value = 287
$ws.Range("H5").Value = 45
$ws.Range("H6").Value = 45

# Fill in the well_depth values (column I)
$ws.Range("I2").Value = 236
$ws.Range("I3").Value = 388
$ws.Range("I4").Value = 639
$ws.Range("I5").Value = 130
$ws.Range("I6").Value = 130

# Match the final cursor/selection position left by the author's edit
$ws.Range("N14").Select() | Out-Null
